# Auto-generated script to apply numeric updates to Siren_Profits workbook
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 107166.664
$ws.Range("I64").Value = 107166.664
$ws.Range("K64").Value = 107166.664
$ws.Range("M64").Value = -106918.664
$ws.Range("H67").Value = 107166.664
$ws.Range("I67").Value = 107166.664
$ws.Range("K67").Value = 107166.664
$ws.Range("M67").Value = -106308.664
$ws.Range("H70").Value = 1684.4286
$ws.Range("I70").Value = 1033.3334
$ws.Range("K70").Value = 3100.0002
$ws.Range("M70").Value = -2830.0002
$ws.Range("H73").Value = 1684.4286
$ws.Range("I73").Value = 1033.3334
$ws.Range("K73").Value = 3100.0002
$ws.Range("M73").Value = -2164.0002
$ws.Range("H95").Value = 55000
$ws.Range("J95").Value = 55000
$ws.Range("L95").Value = 55000
$ws.Range("N95").Value = -60492
$ws.Range("H106").Value = 6325.3
$ws.Range("I106").Value = 6996.857
$ws.Range("J106").Value = 4758.3335
$ws.Range("K106").Value = 6996.857
$ws.Range("L106").Value = 4758.3335
$ws.Range("M106").Value = -6365.857
$ws.Range("N106").Value = -6020.3335
$ws.Range("H121").Value = 2804.3333
$ws.Range("J121").Value = 2804.3333
$ws.Range("L121").Value = 8412.999899999999
$ws.Range("N121").Value = -11906.9999
$ws.Range("H137").Value = 7131.324
$ws.Range("I137").Value = 11365.15
$ws.Range("J137").Value = 2150.353
$ws.Range("K137").Value = 34095.45
$ws.Range("L137").Value = 6451.059
$ws.Range("M137").Value = -31545.45
$ws.Range("N137").Value = -11551.059
$ws.Range("H141").Value = 5702.852
$ws.Range("I141").Value = 5396.913
$ws.Range("K141").Value = 16190.739
$ws.Range("M141").Value = -11010.739

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7785.162
$ws.Range("I61").Value = 10276.682
$ws.Range("K61").Value = 10276.682
$ws.Range("M61").Value = -10064.682
$ws.Range("H122").Value = 2315650.5
$ws.Range("I122").Value = 10692
$ws.Range("J122").Value = 3756249.5
$ws.Range("K122").Value = 32076
$ws.Range("L122").Value = 11268748.5
$ws.Range("M122").Value = -29626
$ws.Range("N122").Value = -11273648.5
$ws.Range("H132").Value = 2728.558
$ws.Range("I132").Value = 2569.5278
$ws.Range("J132").Value = 3546.4285
$ws.Range("K132").Value = 7708.5834
$ws.Range("L132").Value = 10639.2855
$ws.Range("M132").Value = -5178.5834
$ws.Range("N132").Value = -15699.2855
$ws.Range("H136").Value = 7785.162
$ws.Range("I136").Value = 10276.682
$ws.Range("K136").Value = 30830.046
$ws.Range("M136").Value = -28280.046

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 8642177
$ws.Range("J7").Value = 13577207
$ws.Range("L7").Value = 13577207
$ws.Range("N7").Value = -13577433
$ws.Range("H94").Value = 9340.617
$ws.Range("I94").Value = 12443.435
$ws.Range("J94").Value = 2852.9092
$ws.Range("K94").Value = 12443.435
$ws.Range("L94").Value = 2852.9092
$ws.Range("M94").Value = -11992.435
$ws.Range("N94").Value = -3754.9092
$ws.Range("H99").Value = 9358.457
$ws.Range("I99").Value = 9591.034
$ws.Range("K99").Value = 9591.034
$ws.Range("M99").Value = -8093.034
$ws.Range("H105").Value = 56503
$ws.Range("I105").Value = 79343
$ws.Range("J105").Value = 7016.3335
$ws.Range("K105").Value = 79343
$ws.Range("L105").Value = 7016.3335
$ws.Range("M105").Value = -77596
$ws.Range("N105").Value = -10510.3335
$ws.Range("H134").Value = 7278.3184
$ws.Range("I134").Value = 8095.263
$ws.Range("J134").Value = 2104.3333
$ws.Range("K134").Value = 24285.789
$ws.Range("L134").Value = 6312.999899999999
$ws.Range("M134").Value = -21750.789
$ws.Range("N134").Value = -11382.9999
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2812.4138
$ws.Range("I58").Value = 2738.7273
$ws.Range("J58").Value = 3044
$ws.Range("K58").Value = 2738.7273
$ws.Range("L58").Value = 3044
$ws.Range("M58").Value = -2535.7273
$ws.Range("N58").Value = -3450
$ws.Range("H122").Value = 8683.3125
$ws.Range("I122").Value = 8683.3125
$ws.Range("K122").Value = 26049.9375
$ws.Range("M122").Value = -23599.9375
$ws.Range("H134").Value = 11661.077
$ws.Range("I134").Value = 17286.875
$ws.Range("J134").Value = 2659.8
$ws.Range("K134").Value = 51860.625
$ws.Range("L134").Value = 7979.400000000001
$ws.Range("M134").Value = -49325.625
$ws.Range("N134").Value = -13049.4
$ws.Range("H136").Value = 2812.4138
$ws.Range("I136").Value = 2738.7273
$ws.Range("J136").Value = 3044
$ws.Range("K136").Value = 8216.1819
$ws.Range("L136").Value = 9132
$ws.Range("M136").Value = -5666.1819
$ws.Range("N136").Value = -14232
$ws.Range("H139").Value = 75780
$ws.Range("J139").Value = 75780
$ws.Range("L139").Value = 75780
$ws.Range("N139").Value = -86060

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 667818.3
$ws.Range("I5").Value = 475.33334
$ws.Range("J5").Value = 1112713.6
$ws.Range("K5").Value = 1426.00002
$ws.Range("L5").Value = 3338140.8
$ws.Range("M5").Value = -1314.00002
$ws.Range("N5").Value = -3338364.8
$ws.Range("H12").Value = 27.333334
$ws.Range("I12").Value = 64.5
$ws.Range("K12").Value = 193.5
$ws.Range("M12").Value = -20.5
$ws.Range("H74").Value = 17857.143
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -13939
$ws.Range("H77").Value = 17857.143
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -39696
$ws.Range("H92").Value = 403.9091
$ws.Range("I92").Value = 378.2857
$ws.Range("K92").Value = 1134.8571
$ws.Range("M92").Value = 113.1428999999998
$ws.Range("H100").Value = 1200
$ws.Range("I100").Value = 1200
$ws.Range("K100").Value = 3600
$ws.Range("M100").Value = -2789
$ws.Range("H132").Value = 22211.834
$ws.Range("I132").Value = 578.6
$ws.Range("J132").Value = 27904.79
$ws.Range("K132").Value = 5207.400000000001
$ws.Range("L132").Value = 251143.11
$ws.Range("M132").Value = -2677.400000000001
$ws.Range("N132").Value = -256203.11
$ws.Range("H135").Value = 667818.3
$ws.Range("I135").Value = 475.33334
$ws.Range("J135").Value = 1112713.6
$ws.Range("K135").Value = 4278.00006
$ws.Range("L135").Value = 10014422.4
$ws.Range("M135").Value = -1743.00006
$ws.Range("N135").Value = -10019492.4

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3735.2341
$ws.Range("I132").Value = 3827.1538
$ws.Range("K132").Value = 11481.4614
$ws.Range("M132").Value = -8951.4614
$ws.Range("H136").Value = 26045.625
$ws.Range("J136").Value = 26045.625
$ws.Range("L136").Value = 78136.875
$ws.Range("N136").Value = -83236.875

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 54321
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 54321
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 54321
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -55819
$ws.Range("H66").Value = 54321
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 54321
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 162963
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -170451
$ws.Range("H68").Value = 3352.3
$ws.Range("I68").Value = 2099.8572
$ws.Range("J68").Value = 6274.6665
$ws.Range("K68").Value = 2099.8572
$ws.Range("L68").Value = 6274.6665
$ws.Range("M68").Value = -1350.8572
$ws.Range("N68").Value = -7772.6665
$ws.Range("H71").Value = 3352.3
$ws.Range("I71").Value = 2099.8572
$ws.Range("J71").Value = 6274.6665
$ws.Range("K71").Value = 10499.286
$ws.Range("L71").Value = 31373.3325
$ws.Range("M71").Value = -6755.286
$ws.Range("N71").Value = -38861.3325
$ws.Range("H132").Value = 597985.1
$ws.Range("I132").Value = 785628.0600000001
$ws.Range("J132").Value = 3782.5
$ws.Range("K132").Value = 2356884.18
$ws.Range("L132").Value = 11347.5
$ws.Range("M132").Value = -2354354.18
$ws.Range("N132").Value = -16407.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1800
$ws.Range("I96").Value = 900
$ws.Range("K96").Value = 900
$ws.Range("M96").Value = 473
$ws.Range("H122").Value = 3498.537
$ws.Range("I122").Value = 1908.0857
$ws.Range("J122").Value = 6428.316
$ws.Range("K122").Value = 5724.257100000001
$ws.Range("L122").Value = 19284.948
$ws.Range("M122").Value = -3274.257100000001
$ws.Range("N122").Value = -24184.948
$ws.Range("H126").Value = 17169.5
$ws.Range("I126").Value = 20490.904
$ws.Range("J126").Value = 3219.6
$ws.Range("K126").Value = 61472.712
$ws.Range("L126").Value = 9658.799999999999
$ws.Range("M126").Value = -59002.712
$ws.Range("N126").Value = -14598.8
$ws.Range("H132").Value = 7343.1294
$ws.Range("I132").Value = 8718.871999999999
$ws.Range("J132").Value = 3766.2
$ws.Range("K132").Value = 26156.616
$ws.Range("L132").Value = 11298.6
$ws.Range("M132").Value = -23626.616
$ws.Range("N132").Value = -16358.6
$ws.Range("H136").Value = 435280.1
$ws.Range("I136").Value = 444859.5
$ws.Range("K136").Value = 1334578.5
$ws.Range("M136").Value = -1332028.5
